$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Integrated Register")
$ws2 = $wb.Worksheets.Item("Review History")

# --- Sheet1 "Integrated Register": close out risk in row 13 ---
# Status: Open -> closed
$ws1.Range("A13").Value = "closed"
# Date Modified (column E) gets the date of the review (9/23/2015)
$ws1.Range("E13").Value = "9/23/2015"

# Row 6 grew taller (wrapped text reflow)
$ws1.Rows.Item(6).RowHeight = 102.75

# Update the view: scroll/selection moved to C14
$ws1.Activate()
$ws1.Range("C14").Select()

# --- Sheet2 "Review History": log the closure ---
$ws2.Range("A4").Copy()
$ws2.Range("A5").PasteSpecial(-4122)
$ws2.Range("A5").Value = "9/23/2015"
$ws2.Range("B5").Value = "Closed row 13"

$ws2.Activate()
$ws2.Range("A6").Select()

$ws1.Activate()
